# Edit script applying the committed change to
# daten/App_Produkte_mit_EAN.xlsx (sheet "App_Produkte_mit_EAN").
#
# Summary of the change being reproduced:
#   * Three new columns are added after the existing A:D data block:
#       E = "aktiv Cashback" (an "x" marker column, like column C)
#       F = "Cashback single" (numeric cashback amount)
#       G = "Cashback Double" (numeric cashback amount)
#     with header labels in row 1 and per-product values below.
#   * The stray " " marker that lived in cell I10 is moved to J10, and a
#     matching " " marker is added in J9.
#   * A new stand-alone marker row is added at K579 (" ").
#   * The hidden _xlnm._FilterDatabase defined name is expanded from
#     $A$1:$D$1 to cover the sheet's full data range ($A$1:$K$570).
#   * The final selection left on the sheet is the whole header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New/changed cell values (columns E, F, G, J, K) -------------------
# (Row 1 headers are written F, G, E so the newly-created shared-string
#  table entries land in the same order as the target workbook: "Cashback
#  single", "Cashback Double", "aktiv Cashback".)
$ws.Range("F1").Value = "Cashback single"
$ws.Range("G1").Value = "Cashback Double"
$ws.Range("E1").Value = "aktiv Cashback"
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 50
$ws.Range("E3").Value = "x"
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = 50
$ws.Range("E4").Value = "x"
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 50
$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 50
$ws.Range("E6").Value = "x"
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = 50
$ws.Range("E7").Value = "x"
$ws.Range("F7").Value = 50
$ws.Range("G7").Value = 50
$ws.Range("F8").Value = 100
$ws.Range("G8").Value = 200
$ws.Range("E9").Value = "x"
$ws.Range("F9").Value = 100
$ws.Range("G9").Value = 200
$ws.Range("J9").Value = " "
$ws.Range("E10").Value = "x"
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 200
$ws.Range("J10").Value = " "
$ws.Range("F11").Value = 150
$ws.Range("G11").Value = 300
$ws.Range("E12").Value = "x"
$ws.Range("F12").Value = 150
$ws.Range("G12").Value = 300
$ws.Range("E13").Value = "x"
$ws.Range("F13").Value = 150
$ws.Range("G13").Value = 300
$ws.Range("E16").Value = "x"
$ws.Range("F16").Value = 350
$ws.Range("G16").Value = 700
$ws.Range("F17").Value = 150
$ws.Range("G17").Value = 300
$ws.Range("E18").Value = "x"
$ws.Range("F18").Value = 150
$ws.Range("G18").Value = 300
$ws.Range("E19").Value = "x"
$ws.Range("F19").Value = 150
$ws.Range("G19").Value = 300
$ws.Range("F20").Value = 250
$ws.Range("G20").Value = 500
$ws.Range("E21").Value = "x"
$ws.Range("F21").Value = 250
$ws.Range("G21").Value = 500
$ws.Range("E22").Value = "x"
$ws.Range("F22").Value = 250
$ws.Range("G22").Value = 500
$ws.Range("E26").Value = "x"
$ws.Range("F26").Value = 500
$ws.Range("G26").Value = 1000
$ws.Range("F28").Value = 250
$ws.Range("G28").Value = 500
$ws.Range("E29").Value = "x"
$ws.Range("F29").Value = 250
$ws.Range("G29").Value = 500
$ws.Range("E30").Value = "x"
$ws.Range("F30").Value = 250
$ws.Range("G30").Value = 500
$ws.Range("E31").Value = "x"
$ws.Range("F31").Value = 300
$ws.Range("G31").Value = 600
$ws.Range("E32").Value = "x"
$ws.Range("F32").Value = 300
$ws.Range("G32").Value = 600
$ws.Range("F37").Value = 300
$ws.Range("G37").Value = 600
$ws.Range("E38").Value = "x"
$ws.Range("F38").Value = 300
$ws.Range("G38").Value = 600
$ws.Range("E39").Value = "x"
$ws.Range("F39").Value = 500
$ws.Range("G39").Value = 1000
$ws.Range("E40").Value = "x"
$ws.Range("F40").Value = 500
$ws.Range("G40").Value = 1000
$ws.Range("E43").Value = "x"
$ws.Range("F43").Value = 1000
$ws.Range("G43").Value = 2000
$ws.Range("E75").Value = "x"
$ws.Range("F75").Value = 1000
$ws.Range("G75").Value = 2000
$ws.Range("E82").Value = "x"
$ws.Range("E84").Value = "x"
$ws.Range("E85").Value = "x"
$ws.Range("E86").Value = "x"
$ws.Range("E87").Value = "x"
$ws.Range("E92").Value = "x"
$ws.Range("K579").Value = " "

# The stray marker that used to live in I10 moved to J10 above.
$ws.Range("I10").ClearContents()

# --- Expand the hidden AutoFilter range (_FilterDatabase) --------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "App_Produkte_mit_EAN!_FilterDatabase") {
        $n.RefersTo = "=App_Produkte_mit_EAN!`$A`$1:`$K`$570"
    }
}

# --- Leave the sheet selection on the full header row -------------------
$ws.Rows("1:1").Select()
